$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H (copy formatting from the neighboring header cell)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill H2:H21 with a "save" flag: 1 if sum (col G) >= 9, else 0
$lastRow = 21
for ($r = 2; $r -le $lastRow; $r++) {
    $sumVal = $ws.Cells.Item($r, 7).Value2
    if ($sumVal -ge 9) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
